# Add a new "11. References & Resources" section at the end of the
# document body, containing a heading paragraph followed by five
# "Label: <hyperlink>" paragraphs.
#
# We build the new content as a WordprocessingML package fragment and
# insert it with Range.InsertXML so the resulting markup matches the
# target OOXML exactly (plain <w:hyperlink> runs with explicit
# color/underline formatting, no rStyle / w:history noise). The
# placeholder relationship ids below get remapped to the next free
# rIds (rId9..rId13) by the host when the fragment is merged in.

$d = $word.ActiveDocument

$documentXml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:t>11. References &amp; Resources</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">Dataset: </w:t>
  </w:r>
  <w:hyperlink r:id="rIdNewDataset">
    <w:r>
      <w:rPr>
        <w:color w:val="0000FF"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>UCI Machine Learning Repository - Breast Cancer Wisconsin</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">Scikit-learn Documentation: </w:t>
  </w:r>
  <w:hyperlink r:id="rIdNewSklearn">
    <w:r>
      <w:rPr>
        <w:color w:val="0000FF"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>https://scikit-learn.org</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">XGBoost Documentation: </w:t>
  </w:r>
  <w:hyperlink r:id="rIdNewXgboost">
    <w:r>
      <w:rPr>
        <w:color w:val="0000FF"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>https://xgboost.readthedocs.io</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">SMOTE Paper: </w:t>
  </w:r>
  <w:hyperlink r:id="rIdNewSmote">
    <w:r>
      <w:rPr>
        <w:color w:val="0000FF"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>Chawla et al. (2002) - SMOTE</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">GitHub Repository: </w:t>
  </w:r>
  <w:hyperlink r:id="rIdNewGithub">
    <w:r>
      <w:rPr>
        <w:color w:val="0000FF"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>https://github.com/dl1413/LLM-Portfolio</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
</w:body></w:document>
'@

$relsXml = @'
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdNewDataset" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://archive.ics.uci.edu/dataset/17/breast+cancer+wisconsin+diagnostic" TargetMode="External"/>
<Relationship Id="rIdNewSklearn" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://scikit-learn.org" TargetMode="External"/>
<Relationship Id="rIdNewXgboost" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://xgboost.readthedocs.io" TargetMode="External"/>
<Relationship Id="rIdNewSmote" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1613/jair.953" TargetMode="External"/>
<Relationship Id="rIdNewGithub" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/dl1413/LLM-Portfolio" TargetMode="External"/>
</Relationships>
'@

$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
  '<pkg:xmlData>' + $documentXml + '</pkg:xmlData></pkg:part>' +
  '<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">' +
  '<pkg:xmlData>' + $relsXml + '</pkg:xmlData></pkg:part>' +
  '</pkg:package>'

# Collapse a range at the very end of the document body (right after the
# existing "License: MIT" paragraph) and insert the new section there.
$target = $d.Content
$target.Collapse(0)
[void]$target.InsertXML($pkgXml)
